$wb = $excel.ActiveWorkbook

# --- Existing "conect" sheet: the data/content itself is unchanged; only the
# --- view/selection state changes (it stops being the tab-selected sheet and
# --- its lingering selection at K19 is cleared back to the used range A1:E5).
$conect = $wb.Worksheets.Item("conect")

# --- Add the new "test" worksheet right after "conect" ---
$test = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $conect)
$test.Name = "test"

# --- Populate "test" with its header row + matrix values ---
$test.Range("A1").Value = "from"
$test.Range("B1").Value = "quant_x"
$test.Range("C1").Value = "quant_y"
$test.Range("D1").Value = "quant_z"

$test.Range("A2").Value = "demand"
$test.Range("B2").Value = 1
$test.Range("C2").Value = 1
$test.Range("D2").Value = 0

$test.Range("A3").Value = "quant_x"
$test.Range("B3").Value = 0
$test.Range("C3").Value = 0
$test.Range("D3").Value = 0

$test.Range("A4").Value = "quant_y"
$test.Range("B4").Value = 0
$test.Range("C4").Value = 0
$test.Range("D4").Value = 0

$test.Range("A5").Value = "quant_z"
$test.Range("B5").Value = 1
$test.Range("C5").Value = 0
$test.Range("D5").Value = 0

# --- Selection state: "conect" keeps the used range selected ... ---
$conect.Range("A1:E5").Select()

# --- ... while "test" becomes the active/selected tab, with the cursor
# --- left on G7 (outside of the populated range), matching the saved state ---
$test.Activate()
$test.Range("G7").Select()
